$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the weekly-report values between row 3 and row 4 for the
# columns that actually differ (D, M, N, P, R, S). This effectively
# reorders the two weekly records so the earlier-dated record sits in
# row 3 and the later-dated one in row 4.

$cols = @("D", "M", "N", "P", "R", "S")

# Read all the old values first (note: Value is a getter method here,
# must be invoked with parentheses to actually fetch the value).
$row3vals = @{}
$row4vals = @{}
foreach ($col in $cols) {
    $row3vals[$col] = $ws.Range("$col`3").Value()
    $row4vals[$col] = $ws.Range("$col`4").Value()
}

# Now write the swapped values back.
foreach ($col in $cols) {
    $ws.Range("$col`3").Value = $row4vals[$col]
    $ws.Range("$col`4").Value = $row3vals[$col]
}
